# The "estudiantes" sheet still had the placeholder text that was used while
# the allowed-values list for the "Perfil" column was being designed:
#   [Docente, Mentor, Estudiante]
# Now that users can actually log in as "Estudiante", replace that leftover
# placeholder in the Perfil column (H) for every student row with the real
# value "Estudiante".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("estudiantes")

$ws.Range("H2:H4").Value = "Estudiante"
